$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 112-121 (columns A-D numeric values, E unchanged "IMDB reviews")
$data = @{
    112 = @(0.79960707269155218, 0.81243926141885336, 0.81555333998005985, 0.74734982332155475)
    113 = @(0.83622350674373802, 0.83070866141732291, 0.83316880552813433, 0.7859030837004406)
    114 = @(0.8338249754178958, 0.83300589390962665, 0.84189325276938576, 0.78392857142857142)
    115 = @(0.81986368062317427, 0.8176470588235295, 0.8367952522255192, 0.78558875219683644)
    116 = @(0.82793522267206476, 0.81206030150753761, 0.82135523613963035, 0.79491833030852987)
    117 = @(0.83510125361620058, 0.83333333333333337, 0.82189054726368171, 0.78863232682060391)
    118 = @(0.79883381924198249, 0.79377431906614793, 0.80784313725490198, 0.75598935226264419)
    119 = @(0.83011583011583012, 0.82410106899902813, 0.82826300294406296, 0.78540399652476112)
    120 = @(0.82423038728897713, 0.81355932203389825, 0.84084084084084088, 0.79170423805229939)
    121 = @(0.83223992502343003, 0.82521489971346706, 0.82105263157894737, 0.79255319148936165)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
}

# Reset the style for rows 112-121 (columns A-E) back to the default (unstyled) format
$ws.Range("A112:E121").ClearFormats()

# Update the sheet view: remove the frozen/scrolled topLeftCell and selection state
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("A1").Select()
